$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 326 (shifts old 326-330 down to 328-332)
$ws.Rows("326:327").Insert()

# New row 326: Kiwi Hayward, Primera, Región de O'Higgins, fecha 45121
$ws.Range("A326").Value2 = 7
$ws.Range("B326").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C326").Value2 = "Ñuble"
$ws.Range("D326").Value2 = 45121
$ws.Range("E326").Value2 = 16
$ws.Range("F326").Value2 = "Fruta"
$ws.Range("G326").Value2 = 100101
$ws.Range("H326").Value2 = "Berries"
$ws.Range("I326").Value2 = 100101007
$ws.Range("J326").Value2 = "Kiwi"
$ws.Range("K326").Value2 = "Hayward"
$ws.Range("L326").Value2 = "Primera"
$ws.Range("M326").Value2 = 80
$ws.Range("N326").Value2 = 10000
$ws.Range("O326").Value2 = 10000
$ws.Range("P326").Value2 = 10000
$ws.Range("Q326").Value2 = "$/bandeja 18 kilos"
$ws.Range("R326").Value2 = "Región de O'Higgins"
$ws.Range("S326").Value2 = 556
$ws.Range("T326").Value2 = 18

# New row 327: Kiwi Hayward, Segunda, Región de O'Higgins, fecha 45121
$ws.Range("A327").Value2 = 7
$ws.Range("B327").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C327").Value2 = "Ñuble"
$ws.Range("D327").Value2 = 45121
$ws.Range("E327").Value2 = 16
$ws.Range("F327").Value2 = "Fruta"
$ws.Range("G327").Value2 = 100101
$ws.Range("H327").Value2 = "Berries"
$ws.Range("I327").Value2 = 100101007
$ws.Range("J327").Value2 = "Kiwi"
$ws.Range("K327").Value2 = "Hayward"
$ws.Range("L327").Value2 = "Segunda"
$ws.Range("M327").Value2 = 50
$ws.Range("N327").Value2 = 8000
$ws.Range("O327").Value2 = 8000
$ws.Range("P327").Value2 = 8000
$ws.Range("Q327").Value2 = "$/bandeja 18 kilos"
$ws.Range("R327").Value2 = "Región de O'Higgins"
$ws.Range("S327").Value2 = 444
$ws.Range("T327").Value2 = 18
